$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 82, shifting existing rows 82..193 down to 83..194
$ws.Rows("82").Insert()

# Populate the new row 82 with the new record's data
$ws.Range("A82").Value = 4
$ws.Range("B82").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C82").Value = "Los Lagos"
$ws.Range("D82").Value = 44763
$ws.Range("D82").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E82").Value = 10
$ws.Range("F82").Value = 100112009
$ws.Range("G82").Value = "Acelga"
$ws.Range("H82").Value = "Sin especificar"
$ws.Range("I82").Value = "Primera"
$ws.Range("J82").Value = 150
$ws.Range("K82").Value = 1000
$ws.Range("L82").Value = 1000
$ws.Range("M82").Value = 1000
$ws.Range("N82").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O82").Value = "Región de Los Lagos"
$ws.Range("P82").Value = 1000
$ws.Range("Q82").Value = 1
$ws.Range("R82").Value = "Hortaliza"
